$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = '("Endocrine Disruptors"[MeSH Terms] OR "endocrine disrupting chemicals"[Title/Abstract] OR "EDC"[Title/Abstract] OR hormonally active agents[Title/Abstract] OR Endocrine disrupting compounds) AND ("Receptors, Endocrine"[MeSH Terms] OR "receptors"[Title/Abstract] OR "receptor"[Title/Abstract]) AND ("binding"[Title/Abstract] OR "assay"[Title/Abstract] OR "experiment"[Title/Abstract]) NOT review[Publication Type]'
$ws.Range("B10").Value = 1397

$ws.Range("A11").Value = '("Endocrine Disruptors"[MeSH Terms] OR "endocrine disrupting chemicals"[Title/Abstract] OR "EDC"[Title/Abstract] OR "hormonally active agents"[Title/Abstract] OR "Endocrine disrupting compounds"[Title/Abstract]) AND ("Receptors, Endocrine"[MeSH Terms] OR "receptors"[Title/Abstract] OR "receptor"[Title/Abstract]) AND ("binding"[Title/Abstract] OR "assay"[Title/Abstract] OR "experiment"[Title/Abstract]) NOT review[Publication Type] OR (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND ("Bisphenol A"[MeSH Terms] OR "BPA"[Title/Abstract])) OR (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND ("Phthalates"[MeSH Terms] OR "Phthalate"[Title/Abstract])) OR (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND ("Polychlorinated Biphenyls"[MeSH Terms] OR "PCBs"[Title/Abstract])) OR (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND ("Regulation"[MeSH Terms] OR "Risk Assessment"[Title/Abstract]))'
$ws.Range("B11").Value = 5818

$ws.Range("A12").Value = '("Endocrine Disruptors"[MeSH Terms] OR "endocrine disrupting chemicals"[Title/Abstract] OR "EDC"[Title/Abstract] OR "hormonally active agents"[Title/Abstract] OR "Endocrine disrupting compounds"[Title/Abstract]) AND ("Receptors, Endocrine"[MeSH Terms] OR "receptors"[Title/Abstract] OR "receptor"[Title/Abstract]) OR (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND ("Bisphenol A"[MeSH Terms] OR "BPA"[Title/Abstract])) OR (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND ("Phthalates"[MeSH Terms] OR "Phthalate"[Title/Abstract])) OR (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND ("Polychlorinated Biphenyls"[MeSH Terms] OR "PCBs"[Title/Abstract])) OR (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND ("Regulation"[MeSH Terms] OR "Risk Assessment"[Title/Abstract]))'
$ws.Range("B12").Value = 7295

$ws.Range("A13").Value = '("Endocrine Disruptors"[MeSH Terms] OR "endocrine disrupting chemicals"[Title/Abstract] OR "EDC"[Title/Abstract] OR "hormonally active agents"[Title/Abstract] OR "Endocrine disrupting compounds"[Title/Abstract]) AND ("Receptors, Endocrine"[MeSH Terms] OR "receptors"[Title/Abstract] OR "receptor"[Title/Abstract] OR target[Title/Abstract]) AND ("human estrogen receptor alpha"[Title/Abstract] OR  "HERa"[Title/Abstract] OR "Estrogen"[Title/Abstract] OR target[Title/Abstract]) OR "estradiol"[Title/Abstract]AND ("HERb"[Title/Abstract] OR "Human estrogen receptor beta"[Title/Abstract])AND ("binding"[Title/Abstract]AND ("Human andorgen receptor"[Title/Abstract] or "HAR"[Title/Abstract] or "Testosterone"[Title/Abstract])AND ("NR3C1"[Title/Abstract] OR "Human glucocorticoid receptor"[Title/Abstract]) OR  "Dexamethasone"[Title/Abstract]AND ("Human mineralcorticoid receptor"[Title/Abstract],"Aldosterone"[Title/Abstract], "NR3C2"[Tilte/Abstract]AND "Progesterone"[Title/Abstract] OR "NR3C3"[Title/Abstract] OR "Human progesterone receptor"[Title/Abstract] OR "Progesterone"[Title/Abstract]NOT review[Publication Type]OR (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND ("Bisphenol A"[MeSH Terms] OR "BPA"[Title/Abstract])) OR (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND ("Phthalates"[MeSH Terms] OR "Phthalate"[Title/Abstract])) OR (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND ("Polychlorinated Biphenyls"[MeSH Terms] OR "PCBs"[Title/Abstract])) OR (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND ("Regulation"[MeSH Terms] OR "Risk Assessment"[Title/Abstract]))'
$ws.Range("B13").Value = 9999

$ws.Range("A14").Value = '("Endocrine Disruptors"[MeSH Terms] OR "endocrine disrupting chemicals"[Title/Abstract] OR "EDC"[Title/Abstract] OR "hormonally active agents"[Title/Abstract] OR "Endocrine disrupting compounds"[Title/Abstract]) AND (human"[Title/Abstract]")AND ("Receptors, Endocrine"[MeSH Terms] OR "receptors"[Title/Abstract] OR "receptor"[Title/Abstract] OR target[Title/Abstract]) AND ("human estrogen receptor alpha"[Title/Abstract] OR  "HERa"[Title/Abstract] OR "Estrogen"[Title/Abstract] OR target[Title/Abstract]) OR "estradiol"[Title/Abstract]AND ("HERb"[Title/Abstract] OR "Human estrogen receptor beta"[Title/Abstract])AND ("binding"[Title/Abstract]AND ("Human andorgen receptor"[Title/Abstract] or "HAR"[Title/Abstract] or "Testosterone"[Title/Abstract])AND ("NR3C1"[Title/Abstract] OR "Human glucocorticoid receptor"[Title/Abstract]) OR  "Dexamethasone"[Title/Abstract]AND ("Human mineralcorticoid receptor"[Title/Abstract],"Aldosterone"[Title/Abstract], "NR3C2"[Tilte/Abstract]AND "Progesterone"[Title/Abstract] OR "NR3C3"[Title/Abstract] OR "Human progesterone receptor"[Title/Abstract] OR "Progesterone"[Title/Abstract]NOT review[Publication Type]OR (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND ("Bisphenol A"[MeSH Terms] OR "BPA"[Title/Abstract])) OR (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND ("Phthalates"[MeSH Terms] OR "Phthalate"[Title/Abstract])) OR (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND ("Polychlorinated Biphenyls"[MeSH Terms] OR "PCBs"[Title/Abstract])) OR (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND ("Regulation"[MeSH Terms] OR "Risk Assessment"[Title/Abstract]))'
$ws.Range("B14").Value = 9999

$ws.Range("A15").Value = '("Endocrine Disruptors"[MeSH Terms] OR "endocrine disrupting chemicals"[Title/Abstract] OR "EDCs"[Title/Abstract] OR "hormonally active agents"[Title/Abstract] OR "Endocrine disrupting compounds"[Title/Abstract]) AND (human[Title/Abstract]) AND ("Receptors, Endocrine"[MeSH Terms] OR "receptors"[Title/Abstract] OR "receptor"[Title/Abstract] OR target[Title/Abstract]) AND (("human estrogen receptor alpha"[Title/Abstract] OR "HERa"[Title/Abstract] OR "Estrogen"[Title/Abstract] OR target[Title/Abstract]) OR "estradiol"[Title/Abstract] OR ("HERb"[Title/Abstract] OR "Human estrogen receptor beta"[Title/Abstract]) OR ("binding"[Title/Abstract] AND ("Human androgen receptor"[Title/Abstract] OR "HAR"[Title/Abstract] OR "Testosterone"[Title/Abstract])) OR ("NR3C1"[Title/Abstract] OR "Human glucocorticoid receptor"[Title/Abstract]) OR "Dexamethasone"[Title/Abstract] OR ("Human mineralocorticoid receptor"[Title/Abstract] OR "Aldosterone"[Title/Abstract] OR "NR3C2"[Title/Abstract]) OR ("Progesterone"[Title/Abstract] OR "NR3C3"[Title/Abstract] OR "Human progesterone receptor"[Title/Abstract]) ) AND (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND (("Bisphenol A"[MeSH Terms] OR "BPA"[Title/Abstract]) OR ("Phthalates"[MeSH Terms] OR "Phthalate"[Title/Abstract]) OR ("Polychlorinated Biphenyls"[MeSH Terms] OR "PCBs"[Title/Abstract]) OR ("Regulation"[MeSH Terms] OR "Risk Assessment"[Title/Abstract]) ) ) NOT review[Publication Type]'
$ws.Range("B15").Value = 249

$ws.Range("A16").Value = '("Endocrine Disruptors"[MeSH Terms] OR "endocrine disrupting chemicals"[Title/Abstract] OR "EDCs"[Title/Abstract] OR "hormonally active agents"[Title/Abstract] OR "Endocrine disrupting compounds"[Title/Abstract]) AND (human[Title/Abstract]) AND ("Receptors, Endocrine"[MeSH Terms] OR "receptors"[Title/Abstract] OR "receptor"[Title/Abstract] OR target[Title/Abstract]) AND (("human estrogen receptor alpha"[Title/Abstract] OR "HERa"[Title/Abstract] OR "Estrogen"[Title/Abstract] OR target[Title/Abstract]) OR "estradiol"[Title/Abstract] OR ("HERb"[Title/Abstract] OR "Human estrogen receptor beta"[Title/Abstract]) OR ("binding"[Title/Abstract] AND ("Human androgen receptor"[Title/Abstract] OR "HAR"[Title/Abstract] OR "Testosterone"[Title/Abstract])) OR ("NR3C1"[Title/Abstract] OR "Human glucocorticoid receptor"[Title/Abstract]) OR "Dexamethasone"[Title/Abstract] OR ("Human mineralocorticoid receptor"[Title/Abstract] OR "Aldosterone"[Title/Abstract] OR "NR3C2"[Title/Abstract]) OR ("Progesterone"[Title/Abstract] OR "NR3C3"[Title/Abstract] OR "Human progesterone receptor"[Title/Abstract]) ) OR (("Endocrine Disrupting Chemicals"[MeSH Terms] OR "Endocrine Disruptors"[MeSH Terms] OR "EDCs"[Title/Abstract] OR "Endocrine Disruptors"[Title/Abstract]) AND (("Bisphenol A"[MeSH Terms] OR "BPA"[Title/Abstract]) OR ("Phthalates"[MeSH Terms] OR "Phthalate"[Title/Abstract]) OR ("Polychlorinated Biphenyls"[MeSH Terms] OR "PCBs"[Title/Abstract]) ) ) NOT review[Publication Type]'
$ws.Range("B16").Value = 4360
